# Update the "Admin User Input" sheet's name list:
#  - A1 "Vishnu"  -> "Ayaansh" (moved up, takes on the highlighted style from A3)
#  - A3 "Appu"    -> "Mary"    (keeps its existing highlighted style)
#  - A5 "Ayaansh" -> "Leonardo"
# A2 ("Larry") and A4 ("Menon") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A5 before A1 so that the "Ayaansh" shared string already in the
# table is freed up and can be reused/moved naturally when A1 is updated.
$ws.Range("A5").Value = "Leonardo"
$ws.Range("A1").Value = "Ayaansh"
$ws.Range("A3").Value = "Mary"

# Give A1 the same (bold/highlighted) cell style that A3 already has.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Finally, select A5 to match the saved selection state.
$ws.Range("A5").Select() | Out-Null
